$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure "Price" column (D) keeps its original text formatting so numeric-
# looking values like "22.04" are not auto-converted to real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.912.79"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.549.25"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.04"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.770.51"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.550.20"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.914.15"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.61"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.32"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.24"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.415.50"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.969"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  +4.56%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.34"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.683.84"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.70"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +4.83%  "
$ws.Range("E51").Value = "  +0.23%  "

Write-Output "Updated cryptos list"
